# Update the two-digit division worksheet numbers.
#
# The document contains a single 5-column table. Every 4th row is a
# spacer (empty) row, so the "data" rows are rows 1, 5, 9, 13 and 17
# (1-based), each holding 5 division problems such as "80÷3=".
#
# Several of the original expressions (e.g. "67÷7=", "81÷7=") occur
# more than once in the sheet but map to *different* replacements, so a
# single document-wide Find/Replace-All would be ambiguous. We therefore
# target each cell explicitly by (row, column).
#
# Note: a Find.Execute run against a `Cell.Range` / `Table.Cell(...).Range`
# object in this runtime is not confined to that cell the way a literal
# Document.Range(start, end) is, so we rebuild a plain document Range from
# the cell's Start/End before calling Find.Execute (Wrap:=wdFindStop,
# Replace:=wdReplaceOne) to guarantee only the intended cell is touched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-InCell($table, $doc, $row, $col, $old, $new) {
    $cell = $table.Cell($row, $col)
    $start = $cell.Range.Start
    $end = $cell.Range.End
    $rng = $doc.Range($start, $end)
    $found = $rng.Find.Execute($old, $false, $false, $false, $false, `
        $false, $true, 0, $false, $new, 1)
    if (-not $found) {
        Write-Output "WARNING: replacement not found row=$row col=$col old=$old"
    }
}

# Row 1
Replace-InCell $t $d 1 1 "80÷3=" "38÷2="
Replace-InCell $t $d 1 2 "76÷2=" "51÷3="
Replace-InCell $t $d 1 3 "57÷4=" "50÷3="
Replace-InCell $t $d 1 4 "77÷6=" "49÷5="
Replace-InCell $t $d 1 5 "67÷7=" "59÷6="

# Row 5
Replace-InCell $t $d 5 1 "67÷7=" "72÷3="
Replace-InCell $t $d 5 2 "59÷9=" "16÷6="
Replace-InCell $t $d 5 3 "82÷8=" "12÷9="
Replace-InCell $t $d 5 4 "31÷5=" "82÷2="
Replace-InCell $t $d 5 5 "19÷3=" "67÷5="

# Row 9
Replace-InCell $t $d 9 1 "64÷9=" "45÷7="
Replace-InCell $t $d 9 2 "92÷4=" "31÷9="
Replace-InCell $t $d 9 3 "32÷4=" "31÷4="
Replace-InCell $t $d 9 4 "96÷3=" "97÷5="
Replace-InCell $t $d 9 5 "81÷7=" "28÷8="

# Row 13
Replace-InCell $t $d 13 1 "78÷6=" "43÷5="
Replace-InCell $t $d 13 2 "68÷6=" "81÷2="
Replace-InCell $t $d 13 3 "37÷2=" "55÷6="
Replace-InCell $t $d 13 4 "56÷4=" "41÷6="
Replace-InCell $t $d 13 5 "94÷4=" "53÷3="

# Row 17
Replace-InCell $t $d 17 1 "15÷7=" "98÷7="
Replace-InCell $t $d 17 2 "81÷7=" "95÷2="
Replace-InCell $t $d 17 3 "77÷3=" "28÷5="
Replace-InCell $t $d 17 4 "78÷7=" "91÷2="
Replace-InCell $t $d 17 5 "16÷8=" "21÷2="

Write-Output "Done updating division problems."
